# issue #5: add legislator_id, name, date into dataframe
# Target sheet: 股票 (stocks) - append 3 new columns (date, legislator_name,
# legislator_id) to the header row and matching values to the data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("股票")

# --- Header row (row 1): new columns H, I, J ------------------------------
# Copy formatting (bold/border/centered) from the existing header cell G1
# onto the new header cells first, then set their text.
$ws.Range("G1").Copy()
$ws.Range("H1:J1").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("H1").Value = "date"
$ws.Range("I1").Value = "legislator_name"
$ws.Range("J1").Value = "legislator_id"

# --- Data row (row 2): new columns H, I, J --------------------------------
# Copy formatting from the existing data cell G2 onto the new data cells.
$ws.Range("G2").Copy()
$ws.Range("H2:J2").PasteSpecial(-4122)  # xlPasteFormats

# H2 holds a date-looking string ("2012-03-28") that must stay plain text
# (not get auto-converted to a date serial number), so force a text number
# format before assigning it, then restore the plain formatting on top.
$ws.Range("H2").NumberFormat = "@"
$ws.Range("H2").Value = "2012-03-28"
$ws.Range("G2").Copy()
$ws.Range("H2").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("I2").Value = "陳唐山"
$ws.Range("J2").Value = 645

$excel.CutCopyMode = 0

Write-Host "done"
